# Updated cryptos list on Sun May 12 21:50:37 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price (D) and Volume (E) columns to text format so
# numeric-looking strings (e.g. "0.999", "1.00", "61.254.06") are written as
# literal text instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "61.254.06"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.930.57"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "595.56"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "144.03"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.501"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "6.97"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "33.24"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "3.412.30"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "61.207.87"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "2.925.76"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "6.65"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "431.99"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "13.51"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").Value = "7.05"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "81.57"
$ws.Range("D24").Value = "10.90"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").Value = "11.73"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").Value = "2.59"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").Value = "6.90"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("D31").Value = "26.60"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "0.0₃0874"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "5.61"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "2.96"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").Value = "42.16"
$ws.Range("E41").Value = "  +4.89%  "
$ws.Range("D42").Value = "0.279"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "2.690.78"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "133.58"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").Value = "363.39"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D48").Value = "23.54"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "2.00"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -0.94%  "

# Restore the default cell style so no formatting/style metadata changes leak
# into the saved workbook (matches the original, unstyled data cells).
$dataRange.Style = "Normal"
